# Update cryptos list with latest price/volume figures (2023-02-25 run)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "23.056.76"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -0.11%  "

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.591.09"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -0.45%  "

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.002"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.13%  "

# Row 5
$ws.Range("E5").Value = "  -0.03%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "302.10"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.12%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.3773"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.03%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3604"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -1.38%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "50.92"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +6.45%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.002"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -0.16%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1.231"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -3.55%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.08073"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.03%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "22.11"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -3.63%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.499"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -2.04%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.284"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -4.69%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.00001229"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -2.85%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "1.590.93"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.94%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "92.67"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +1.15%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06815"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.33%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "18.01"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -2.22%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.473"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -1.79%  "

# Row 22
$ws.Range("E22").Value = "  +0.08%  "

# Row 23
$ws.Range("E23").Value = "  -1.23%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "23.061.74"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.14%  "

# Row 25
$ws.Range("E25").Value = "  +0.12%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.826"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -2.34%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "20.95"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.66%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "148.60"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -1.59%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "5.216"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.55%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "133.54"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +1.08%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "2.353"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -3.80%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "6.575"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -7.48%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.767.24"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -0.23%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.9482"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -3.46%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.07424"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -3.94%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "10.14"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +0.88%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.02684"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -3.43%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.08788"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -0.97%  "

# Row 39
$ws.Range("B39").Value = "InternetComputer(DFINITY)"
$ws.Range("C39").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "6.083"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -3.51%  "

# Row 40
$ws.Range("B40").Value = "Algorand"
$ws.Range("C40").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.2485"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -2.26%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.348"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -3.52%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.6945"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -2.90%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "12.13"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -5.09%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "14.94"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -5.71%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.6473"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -2.52%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "4.010"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +1.09%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.260"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -2.21%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "131.83"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.21%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.07897"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.95%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.203"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +2.60%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.215"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +4.06%  "
